$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price column so numeric-looking strings
# (e.g. "1.000") are preserved as literal text instead of being
# coerced into numbers by Excel's automatic type detection.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '30.324.61'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '1.870.44'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '244.97'
$ws.Range('E5').Value = '  +4.34%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = '0.4720'
$ws.Range('E7').Value = '  +0.44%  '
$ws.Range('D8').Value = '0.2872'
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('D9').Value = '0.06470'
$ws.Range('E9').Value = '  -1.60%  '
$ws.Range('D10').Value = '21.09'
$ws.Range('E10').Value = '  -1.29%  '
$ws.Range('D11').Value = '0.07773'
$ws.Range('E11').Value = '  -0.80%  '
$ws.Range('D12').Value = '1.867.66'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').Value = '95.23'
$ws.Range('E13').Value = '  -1.84%  '
$ws.Range('D14').Value = '0.7156'
$ws.Range('E14').Value = '  +2.52%  '
$ws.Range('D15').Value = '5.107'
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('D16').Value = '276.53'
$ws.Range('E16').Value = '  +2.70%  '
$ws.Range('D17').Value = '30.308.98'
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').Value = '13.36'
$ws.Range('E18').Value = '  -3.28%  '
$ws.Range('D19').Value = '0.000007555'
$ws.Range('E19').Value = '  -1.11%  '
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').Value = '2.119.00'
$ws.Range('E21').Value = '  +0.36%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').Value = '5.228'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').Value = '6.159'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').Value = '9.266'
$ws.Range('E25').Value = '  -1.86%  '
$ws.Range('D26').Value = '165.54'
$ws.Range('E26').Value = '  -1.09%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').Value = '1.907'
$ws.Range('E28').Value = '  -1.75%  '
$ws.Range('D29').Value = '1.385'
$ws.Range('E29').Value = '  +1.32%  '
$ws.Range('D30').Value = '0.09908'
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('D31').Value = '1.516'
$ws.Range('E31').Value = '  +3.89%  '
$ws.Range('D32').Value = '4.275'
$ws.Range('E32').Value = '  -1.84%  '
$ws.Range('D33').Value = '4.028'
$ws.Range('E33').Value = '  -0.54%  '
$ws.Range('D34').Value = '0.04766'
$ws.Range('E34').Value = '  +0.87%  '
$ws.Range('D35').Value = '1.123'
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('D36').Value = '0.6935'
$ws.Range('E36').Value = '  -1.38%  '
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').Value = '0.01851'
$ws.Range('E38').Value = '  -1.18%  '
$ws.Range('D39').Value = '2.739'
$ws.Range('E39').Value = '  -0.97%  '
$ws.Range('D40').Value = '6.354'
$ws.Range('E40').Value = '  +0.68%  '
$ws.Range('D41').Value = '70.30'
$ws.Range('E41').Value = '  -3.49%  '
$ws.Range('D42').Value = '1.914'
$ws.Range('E42').Value = '  -1.80%  '
$ws.Range('D43').Value = '0.8418'
$ws.Range('E43').Value = '  +0.89%  '
$ws.Range('D44').Value = '0.9999'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').Value = '0.4112'
$ws.Range('E45').Value = '  -1.51%  '
$ws.Range('D46').Value = '102.00'
$ws.Range('E46').Value = '  -1.12%  '
$ws.Range('D47').Value = '9.313'
$ws.Range('E47').Value = '  +2.16%  '
$ws.Range('D48').Value = '7.102'
$ws.Range('E48').Value = '  -0.10%  '
$ws.Range('D49').Value = '35.26'
$ws.Range('E49').Value = '  +2.27%  '
$ws.Range('D50').Value = '920.10'
$ws.Range('E50').Value = '  -5.20%  '
$ws.Range('E51').Value = '  -1.92%  '

# Restore the default (Normal) style on the Price column so no stray
# text-format style lingers on cells after the text has been written.
$ws.Range("D2:D51").Style = "Normal"
